$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Cells.Item(2, 4)
$cell.NumberFormat = "@"
$cell.Value = "23.306.24"
$cell.Style = "Normal"
$ws.Cells.Item(2, 5).Value = "  -0.60%  "

$cell = $ws.Cells.Item(3, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.622.12"
$cell.Style = "Normal"
$ws.Cells.Item(3, 5).Value = "  -1.15%  "

$cell = $ws.Cells.Item(4, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.002"
$cell.Style = "Normal"
$ws.Cells.Item(4, 5).Value = "  +0.19%  "

$ws.Cells.Item(5, 5).Value = "  +0.10%  "

$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = "302.40"
$cell.Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -0.92%  "

$cell = $ws.Cells.Item(7, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.3748"
$cell.Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  +0.44%  "

$ws.Cells.Item(8, 2).Value = "OKB"
$ws.Cells.Item(8, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$cell = $ws.Cells.Item(8, 4)
$cell.NumberFormat = "@"
$cell.Value = "51.40"
$cell.Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  -1.87%  "

$ws.Cells.Item(9, 2).Value = "Cardano"
$ws.Cells.Item(9, 3).Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.3615"
$cell.Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  -0.11%  "

$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.08138"
$cell.Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  +0.21%  "

$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.219"
$cell.Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  -2.72%  "

$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.002"
$cell.Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  +0.20%  "

$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = "@"
$cell.Value = "22.26"
$cell.Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  -2.51%  "

$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = "6.455"
$cell.Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  -2.18%  "

$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.00001234"
$cell.Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  -2.75%  "

$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = "@"
$cell.Value = "7.251"
$cell.Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  -0.43%  "

$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.618.56"
$cell.Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  -0.76%  "

$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = "@"
$cell.Value = "93.88"
$cell.Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  -0.49%  "

$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.06926"
$cell.Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  +0.92%  "

$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = "17.50"
$cell.Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  -3.47%  "

$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = "6.517"
$cell.Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  +0.08%  "

$ws.Cells.Item(22, 5).Value = "  +0.16%  "

$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = "12.47"
$cell.Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  -2.07%  "

$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = "23.311.29"
$cell.Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  -0.56%  "

$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.469"
$cell.Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  +2.61%  "

$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.078"
$cell.Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  +1.37%  "

$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = "21.09"
$cell.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  -0.62%  "

$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = "@"
$cell.Value = "150.54"
$cell.Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  -0.79%  "

$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = "@"
$cell.Value = "5.267"
$cell.Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  -0.60%  "

$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = "@"
$cell.Value = "132.54"
$cell.Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  -2.14%  "

$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.797.84"
$cell.Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  -0.77%  "

$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = "@"
$cell.Value = "6.715"
$cell.Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  -0.60%  "

$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.166"
$cell.Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  -5.46%  "

$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.063"
$cell.Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  +11.51%  "

$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = "@"
$cell.Value = "11.20"
$cell.Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  +8.30%  "

$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.02743"
$cell.Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  -3.35%  "

$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.08772"
$cell.Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  -0.12%  "

$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.2475"
$cell.Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  -1.76%  "

$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.07081"
$cell.Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  -2.07%  "

$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = "5.971"
$cell.Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  -1.37%  "

$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.6958"
$cell.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  -1.31%  "

$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.332"
$cell.Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  -3.19%  "

$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = "15.97"
$cell.Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  -0.85%  "

$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = "11.98"
$cell.Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  -3.83%  "

$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.6444"
$cell.Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  -1.06%  "

$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.001"
$cell.Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  +0.06%  "

$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.953"
$cell.Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  -1.42%  "

$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.258"
$cell.Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  -2.88%  "

$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.07953"
$cell.Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  -0.19%  "

$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = "125.65"
$cell.Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  -2.07%  "

$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.181"
$cell.Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  -1.35%  "
